{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Change being applied (per the supplied OOXML diff):\n//  1) Remove the paragraph that holds \"CAROLINE WAMBUI \" + \"AS263/0858/2019\"\n//     (this paragraph also carried the empty `_GoBack` bookmark around the\n//     first run). Deleting the whole paragraph merges it away, leaving the\n//     single blank paragraph that used to follow it.\n//  2) Re-create the (still empty) `_GoBack` bookmark, now anchored at the\n//     end of the \"X =   [, 1] [, 2] [, 3]     \" paragraph, i.e. right after\n//     its last run.\n//\n// (The rest of the raw diff - around \"print(Transpose.X)\" / \"Output:\" - is\n// a by-product of the diffing tool re-aligning unrelated, unchanged lines\n// after the edit above shifted line numbers; that region's content is\n// identical before and after, so nothing needs to change there.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet signaturePara = null;\nlet matrixPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (signaturePara === null && text.indexOf(\"CAROLINE WAMBUI\") !== -1) {\n    signaturePara = para;\n  }\n  if (\n    matrixPara === null &&\n    text.indexOf(\"X =\") !== -1 &&\n    text.indexOf(\"[, 1] [, 2] [, 3]\") !== -1\n  ) {\n    matrixPara = para;\n  }\n  if (signaturePara && matrixPara) break;\n}\n\nif (signaturePara) {\n  signaturePara.delete();\n}\n\nif (matrixPara) {\n  const endRange = matrixPara.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Change being applied (per the supplied OOXML diff):\n#  1) Remove the paragraph that holds \"CAROLINE WAMBUI \" + \"AS263/0858/2019\"\n#     (this paragraph also carried the empty `_GoBack` bookmark around the\n#     first run). Deleting the whole paragraph (its Range, pilcrow included)\n#     merges it away, leaving the single blank paragraph that used to follow\n#     it.\n#  2) Re-create the (still empty) `_GoBack` bookmark, now anchored at the\n#     end of the \"X =   [, 1] [, 2] [, 3]     \" paragraph, i.e. right after\n#     its last run (before the paragraph mark).\n#\n# (The rest of the raw diff - around \"print(Transpose.X)\" / \"Output:\" - is\n# a by-product of the diffing tool re-aligning unrelated, unchanged lines\n# after the edit above shifted line numbers; that region's content is\n# identical before and after, so nothing needs to change there.)\n\n$d = $word.ActiveDocument\n\n# --- locate the two paragraphs we care about by their text -----------------\n$signatureIndex = -1\n$matrixIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($signatureIndex -eq -1 -and $t.IndexOf(\"CAROLINE WAMBUI\") -ge 0) {\n        $signatureIndex = $i\n    }\n    if ($matrixIndex -eq -1 -and $t.IndexOf(\"X =\") -ge 0 -and $t.IndexOf(\"[, 1] [, 2] [, 3]\") -ge 0) {\n        $matrixIndex = $i\n    }\n    if ($signatureIndex -ne -1 -and $matrixIndex -ne -1) {\n        break\n    }\n}\n\n# --- step 1: delete the whole \"CAROLINE WAMBUI ... AS263/0858/2019\" paragraph\nif ($signatureIndex -ne -1) {\n    $d.Paragraphs.Item($signatureIndex).Range.Delete()\n}\n\n# --- step 2: re-add an empty `_GoBack` bookmark at the end of the matrix paragraph\nif ($matrixIndex -ne -1) {\n    $matrixPara = $d.Paragraphs.Item($matrixIndex)\n    $paraRange = $matrixPara.Range\n    $endOfText = $paraRange.End - 1   # position right before the paragraph mark\n\n    # A zero-length Range placed exactly at (paragraph-end - 1) is mis-resolved\n    # by this host when handed straight to Bookmarks.Add, so we briefly insert\n    # a one-character marker there, wrap the bookmark around that marker\n    # (a non-empty range, which anchors correctly), then delete the marker -\n    # the bookmark collapses to an empty range at the right spot and survives.\n    $insertionPoint = $d.Range($endOfText, $endOfText)\n    $insertionPoint.InsertAfter([char]1)\n\n    $matrixPara2 = $d.Paragraphs.Item($matrixIndex)\n    $paraRange2 = $matrixPara2.Range\n    $markerRange = $d.Range($paraRange2.End - 2, $paraRange2.End - 1)\n\n    $d.Bookmarks.Add(\"_GoBack\", $markerRange)\n    $d.Bookmarks.Item(\"_GoBack\").Range.Delete()\n}\n"}
